$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "canonical SMILES" column (D) ---
# Copy the formatting (styles) from column C onto column D first, so the new
# column inherits the same borders/fills/fonts/alignment as the rest of the
# table, then overwrite the values with the non-isomeric SMILES.
$ws.Range("C2:C28").Copy()
$ws.Range("D2:D28").PasteSpecial(-4122)

# Header
$ws.Range("D2").Value = "canonical SMILES"

# Body rows: "canonical SMILES" is the isomeric SMILES with stereo-bond
# markers ("/" and "\") stripped out.
for ($r = 3; $r -le 28; $r++) {
    $iso = $ws.Cells.Item($r, 3).Value()
    if ($iso -ne $null) {
        $plain = $iso.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $plain
    }
}

# Column D width to roughly match the authored width (47.28515625 chars).
$ws.Columns.Item(4).ColumnWidth = 46.42
